# Update the cryptos worksheet with the latest scraped prices / 1h-volume
# deltas (and restore the two rows whose coins were reordered upstream).
#
# Several "Price" cells contain plain numeric-looking text (e.g. "1.001")
# that Excel's COM layer would otherwise silently coerce to a Double on
# assignment (dropping the significant trailing zero / changing type from
# string to number). Force those cells to text first so the literal string
# is preserved exactly as scraped, then drop the now-unneeded explicit
# "@" number format so the cell's style stays untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

$ws.Range("D2").Value = "30.712.25"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "1.889.48"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  +2.07%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  +0.13%  "
Set-TextValue $ws.Range("D7") "0.4938"
$ws.Range("E7").Value = "  +0.20%  "
Set-TextValue $ws.Range("D8") "0.2960"
$ws.Range("E8").Value = "  +1.38%  "
Set-TextValue $ws.Range("D9") "0.06818"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").Value = "1.889.92"
$ws.Range("E10").Value = "  +0.71%  "
Set-TextValue $ws.Range("D11") "17.22"
$ws.Range("E11").Value = "  +2.87%  "
Set-TextValue $ws.Range("D12") "0.07244"
$ws.Range("E12").Value = "  -0.03%  "
Set-TextValue $ws.Range("D13") "91.62"
$ws.Range("E13").Value = "  +6.25%  "
Set-TextValue $ws.Range("D14") "5.078"
$ws.Range("E14").Value = "  +3.93%  "
Set-TextValue $ws.Range("D15") "0.6788"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "30.681.87"
$ws.Range("E16").Value = "  +2.46%  "
Set-TextValue $ws.Range("D17") "0.000007989"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("E18").Value = "  +0.09%  "
Set-TextValue $ws.Range("D19") "13.25"
$ws.Range("E19").Value = "  +3.82%  "
$ws.Range("D20").Value = "2.132.13"
$ws.Range("E20").Value = "  +0.56%  "
Set-TextValue $ws.Range("D21") "1.003"
$ws.Range("E21").Value = "  +0.36%  "
Set-TextValue $ws.Range("D22") "4.832"
$ws.Range("E22").Value = "  +1.22%  "
Set-TextValue $ws.Range("D23") "187.42"
$ws.Range("E23").Value = "  +32.23%  "
Set-TextValue $ws.Range("D24") "6.054"
$ws.Range("E24").Value = "  +4.91%  "
Set-TextValue $ws.Range("D25") "9.354"
$ws.Range("E25").Value = "  +3.33%  "
Set-TextValue $ws.Range("D26") "156.19"
$ws.Range("E26").Value = "  +4.33%  "
Set-TextValue $ws.Range("D27") "19.16"
$ws.Range("E27").Value = "  +12.53%  "
Set-TextValue $ws.Range("D28") "1.910"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  +0.51%  "
Set-TextValue $ws.Range("D30") "4.309"
Set-TextValue $ws.Range("D31") "0.09009"
$ws.Range("E31").Value = "  +3.13%  "
Set-TextValue $ws.Range("D32") "4.012"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("E33").Value = "  +2.52%  "
Set-TextValue $ws.Range("D34") "0.7457"
$ws.Range("E34").Value = "  +4.76%  "
Set-TextValue $ws.Range("D35") "1.116"
$ws.Range("E35").Value = "  +0.05%  "
Set-TextValue $ws.Range("D36") "2.736"
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("E37").Value = "  +3.00%  "
Set-TextValue $ws.Range("D38") "2.667"
$ws.Range("E38").Value = "  -0.73%  "
Set-TextValue $ws.Range("D39") "2.160"
$ws.Range("E39").Value = "  -0.70%  "
Set-TextValue $ws.Range("D40") "0.9393"
$ws.Range("E40").Value = "  +0.97%  "
Set-TextValue $ws.Range("D41") "0.4427"
$ws.Range("E41").Value = "  +4.38%  "
Set-TextValue $ws.Range("D42") "105.46"
$ws.Range("E42").Value = "  +2.64%  "
Set-TextValue $ws.Range("D43") "5.774"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("E44").Value = "  +0.23%  "
Set-TextValue $ws.Range("D45") "7.630"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("E46").Value = "  +5.75%  "
Set-TextValue $ws.Range("D47") "0.05845"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "1.434"
$ws.Range("E48").Value = "  +7.53%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "8.641"
$ws.Range("E49").Value = "  +4.30%  "
Set-TextValue $ws.Range("D50") "0.3939"
$ws.Range("E50").Value = "  +4.22%  "
Set-TextValue $ws.Range("D51") "33.53"
$ws.Range("E51").Value = "  +2.99%  "
